$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3065615121153629
$ws.Range("C2").Value = 1.723407689173405
$ws.Range("D2").Value = 11.71420312276167
$ws.Range("E2").Value = 3.4226018060478
$ws.Range("F2").Value = 3.489063932081051
$ws.Range("G2").Value = 22
